# Update "想去人数" (interest count, column F) figures across all four
# sheets to the freshly scraped values (gh-pages data regen at 456a3b4).
# Only column F values change; everything else in the workbook is untouched.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 765
$ws.Range("F3").Value = 965
$ws.Range("F4").Value = 750
$ws.Range("F5").Value = 851
$ws.Range("F6").Value = 419
$ws.Range("F7").Value = 636
$ws.Range("F8").Value = 140
$ws.Range("F9").Value = 1237
$ws.Range("F10").Value = 663
$ws.Range("F11").Value = 394
$ws.Range("F12").Value = 521
$ws.Range("F14").Value = 20
$ws.Range("F15").Value = 673
$ws.Range("F16").Value = 6
$ws.Range("F17").Value = 375
$ws.Range("F19").Value = 85
$ws.Range("F20").Value = 564
$ws.Range("F21").Value = 105
$ws.Range("F22").Value = 601
$ws.Range("F24").Value = 845

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 52
$ws.Range("F13").Value = 102

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 371

# 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 371
$ws.Range("F6").Value = 765
$ws.Range("F7").Value = 965
$ws.Range("F8").Value = 750
$ws.Range("F9").Value = 851
$ws.Range("F10").Value = 419
$ws.Range("F11").Value = 636
$ws.Range("F12").Value = 140
$ws.Range("F13").Value = 1237
$ws.Range("F14").Value = 663
$ws.Range("F17").Value = 394
$ws.Range("F18").Value = 521
$ws.Range("F21").Value = 20
$ws.Range("F22").Value = 673
$ws.Range("F24").Value = 6
$ws.Range("F25").Value = 375
$ws.Range("F27").Value = 85
$ws.Range("F29").Value = 52
$ws.Range("F30").Value = 564
$ws.Range("F33").Value = 102
$ws.Range("F34").Value = 102
$ws.Range("F35").Value = 105
$ws.Range("F36").Value = 601
$ws.Range("F38").Value = 845
